# Add new column 'Servised by' to Card22 by admin
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card22")

# New header cell O1: same text style as the existing header N1 ("Correction")
$ws.Cells.Item(1, 15).Value = "Servised by"
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)   # xlPasteFormats

# Fill in the previously-blank N column (rows 2-12) with "nan" to match
# the rest of the sheet's placeholder convention, and create the new
# (blank) O column cells alongside them.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"

    # Create O<r> as a genuine blank text cell (matching the type/style of
    # the sheet's other blank cells, e.g. M<r>) rather than leaving it absent.
    $ws.Cells.Item($r, 15).Value = "'"
    $ws.Range("M$r").Copy()
    $ws.Range("O$r").PasteSpecial(-4122)   # xlPasteFormats
}

$excel.CutCopyMode = 0
